$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.990.18"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "3.400.49"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").Value = "3.400.34"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "3.978.11"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "3.388.86"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "61.041.17"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.559"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.33%  "
$ws.Range("D27").Value = "3.534.29"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").Value = "3.429.20"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("E40").Value = "  -5.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0773"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -4.39%  "
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "2.491.27"
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  +0.44%  "
